$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used after "Jos Buttler" in the source data (matches
# existing rows 2-6 in the sheet).
$nbsp = [char]0x00A0
$batsman = "Jos Buttler" + $nbsp

# New rows 7-11 duplicate existing rows 5, 3, 4, 6, 2 (in that order).
$newRows = @(
    @(" Abu Dhabi", " October 30 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", $batsman, "22", "11", "1", "2", "200.00"),
    @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", $batsman, "9", "12", "0", "0", "75.00"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Rajasthan Royals", "Kolkata Knight Riders", $batsman, "35", "22", "4", "1", "159.09"),
    @(" Abu Dhabi", " October 19 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Rajasthan Royals", "Chennai Super Kings", $batsman, "70", "48", "7", "2", "145.83"),
    @(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", $batsman, "24", "25", "1", "1", "96.00")
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        # Leading apostrophe forces Excel to store the value as text, even
        # for the numeric-looking columns (G:K), matching the workbook's
        # existing "number stored as text" cells.
        $ws.Cells.Item($rowNum, $col + 1).Value = "'" + $rowValues[$col]
    }
}
